$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @{
  2 = @(1.02, 1.028455167873632, 1.037090373158348, 1.028400354951506, 1.03922080121117, 1.036307416417705, 1.033607391713031, 1.039882250462479, 1.031217281667805, 1.042006606883858, 1.01514464212819)
  3 = @(1.02, 1.029567479636824, 1.037697629840338, 1.029350129067517, 1.04052663275196, 1.036562816819274, 1.034359485292177, 1.040299779767379, 1.031974600653153, 1.043121305582433, 1.015399913991212)
  4 = @(1.02, 1.03028708488381, 1.038090497351827, 1.02996491696645, 1.041371745731729, 1.036726829572157, 1.034845480418887, 1.040569186866336, 1.032464240192621, 1.043842201309539, 1.015564715727895)
  5 = @(1.02, 1.030589575823875, 1.038255641605806, 1.030223426899438, 1.041727069712912, 1.036795481641341, 1.035049635355809, 1.040682262651516, 1.03266999011904, 1.044145174364458, 1.015633908370531)
  6 = @(1.02, 1.030640363593891, 1.038283369014195, 1.03026683497731, 1.041786732486002, 1.036806991097897, 1.035083904652342, 1.040701237830931, 1.032704530904517, 1.044196039553574, 1.015645520843405)
  7 = @(1.02, 1.030291126903303, 1.038092704086382, 1.029968370979802, 1.041376493434512, 1.03672774807894, 1.034848208965027, 1.04057069851001, 1.032466989804193, 1.043846250007911, 1.015565640636963)
  8 = @(1.02, 1.028831107469255, 1.037295611955465, 1.028721290080843, 1.039662082652586, 1.036393988650064, 1.033861702171144, 1.04002351400104, 1.031473303490707, 1.042383406000887, 1.015230990388)
  9 = @(1.02, 1.026257280993002, 1.035890552250591, 1.026525449069736, 1.036642153972861, 1.035796301615852, 1.032118274775497, 1.039053487102259, 1.029719247011357, 1.039802638383736, 1.014638409402911)
  10 = @(1.02, 1.024540590164494, 1.0349535753413, 1.025062658469171, 1.03462947186796, 1.035391415583718, 1.030952541195019, 1.038402916838608, 1.028547795488957, 1.038079974335114, 1.014241410786507)
  11 = @(1.02, 1.023797034963685, 1.034547800660625, 1.024429509763243, 1.033758073188582, 1.035214570313775, 1.030446939061507, 1.038120296122543, 1.028040043803816, 1.037333507886427, 1.014069043270219)
  12 = @(1.02, 1.023520811255833, 1.034397070062474, 1.02419436698261, 1.033434410306286, 1.035148652590916, 1.030259010110397, 1.038015180405414, 1.027851365608928, 1.037056153536776, 1.01400494834617)
  13 = @(1.02, 1.023580063752729, 1.03442940262793, 1.024244804250056, 1.033503836528182, 1.035162802554059, 1.030299327214863, 1.038037734313788, 1.027891841199898, 1.037115650797621, 1.014018700093239)
  14 = @(1.02, 1.023774202932956, 1.034535341377042, 1.024410072039675, 1.033731318855346, 1.035209126216923, 1.030431407367061, 1.038111610033663, 1.028024449172126, 1.037310583404443, 1.014063746591821)
  15 = @(1.02, 1.023893813908554, 1.03460061272578, 1.024511903842646, 1.033871479922982, 1.035237637338166, 1.030512769612286, 1.038157109053355, 1.028106143148468, 1.037430676692703, 1.014091491951683)
  16 = @(1.02, 1.024589932768183, 1.034980504129557, 1.025104683663449, 1.034687305744705, 1.035403120041561, 1.030986078746902, 1.03842165410149, 1.028581482589455, 1.03812950326483, 1.014252840450142)
  17 = @(1.02, 1.025026530555359, 1.035218785161028, 1.025476585092589, 1.035199078310024, 1.03550651408593, 1.03128274977534, 1.038587350289918, 1.028879514665964, 1.038567712233196, 1.014353925623426)
  18 = @(1.02, 1.025281169855257, 1.035357764852969, 1.025693532997018, 1.035497596981925, 1.03556667476194, 1.031455712766662, 1.038683909332089, 1.02905330301112, 1.038823260000918, 1.014412842089106)
  19 = @(1.02, 1.025367991797282, 1.035405152360947, 1.025767510731092, 1.035599385974437, 1.035587163012912, 1.031514675099027, 1.03871681840756, 1.029112552101947, 1.038910386411292, 1.014432923486645)
  20 = @(1.02, 1.024979689908015, 1.035193220434882, 1.025436681101092, 1.035144168900098, 1.035495436122286, 1.031250928083036, 1.038569581833393, 1.028847543689039, 1.038520701962977, 1.014343084776617)
  21 = @(1.02, 1.023717034740396, 1.034504145265028, 1.024361403782798, 1.033664330630673, 1.035195491389431, 1.030392516533879, 1.038089859276656, 1.02798540154975, 1.037253182931581, 1.014050483453175)
  22 = @(1.02, 1.022922954642578, 1.034070851545651, 1.023685547092628, 1.032733973768454, 1.035005576270074, 1.029852070826335, 1.037787441567342, 1.027442894982645, 1.036455759404606, 1.013866108571726)
  23 = @(1.02, 1.023343930906558, 1.034300552719033, 1.02404381154252, 1.033227167168665, 1.035106379780276, 1.030138640636318, 1.037947834358653, 1.027730530251231, 1.03687853528859, 1.013963887613567)
  24 = @(1.02, 1.025000855244133, 1.035204772051557, 1.025454711922509, 1.035168980066123, 1.035500442232321, 1.03126530718356, 1.038577610908232, 1.028861990147632, 1.038541944040993, 1.014347983427014)
  25 = @(1.02, 1.026922813367461, 1.036253845406673, 1.027092930444492, 1.037422762431828, 1.035951951000466, 1.032569597776853, 1.039304949065932, 1.030173077370371, 1.040470200552347, 1.014791948152948)
}

$cols = @(2,3,4,5,6,9,10,11,12,13,14)

foreach ($r in $data.Keys) {
  $rowVals = $data[$r]
  for ($i = 0; $i -lt $cols.Count; $i++) {
    $ws.Cells.Item([int]$r, $cols[$i]).Value = $rowVals[$i]
  }
}